$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "98×58=5684" },
    @{ Row = 1;  Col = 2; Text = "81×24=1944" },
    @{ Row = 1;  Col = 3; Text = "47×24=1128" },
    @{ Row = 1;  Col = 4; Text = "39×67=2613" },
    @{ Row = 1;  Col = 5; Text = "39×25=975" },
    @{ Row = 5;  Col = 1; Text = "86×98=8428" },
    @{ Row = 5;  Col = 2; Text = "35×11=385" },
    @{ Row = 5;  Col = 3; Text = "75×84=6300" },
    @{ Row = 5;  Col = 4; Text = "19×83=1577" },
    @{ Row = 5;  Col = 5; Text = "94×82=7708" },
    @{ Row = 10; Col = 1; Text = "71×40=2840" },
    @{ Row = 10; Col = 2; Text = "29×60=1740" },
    @{ Row = 10; Col = 3; Text = "87×26=2262" },
    @{ Row = 10; Col = 4; Text = "61×96=5856" },
    @{ Row = 10; Col = 5; Text = "86×53=4558" },
    @{ Row = 15; Col = 1; Text = "96×63=6048" },
    @{ Row = 15; Col = 2; Text = "18×39=702" },
    @{ Row = 15; Col = 3; Text = "32×87=2784" },
    @{ Row = 15; Col = 4; Text = "90×88=7920" },
    @{ Row = 15; Col = 5; Text = "28×26=728" },
    @{ Row = 20; Col = 1; Text = "52×57=2964" },
    @{ Row = 20; Col = 2; Text = "71×74=5254" },
    @{ Row = 20; Col = 3; Text = "20×68=1360" },
    @{ Row = 20; Col = 4; Text = "86×53=4558" },
    @{ Row = 20; Col = 5; Text = "70×95=6650" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}

Write-Host "Updated" $updates.Count "cells"
